$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new formula in D4: ISBLANK(C4) -> TRUE (C4 is empty)
$ws.Range("D4").Formula = "=ISBLANK(C4)"

# Add new formula in C5 referencing an undefined name -> #NAME? error
$ws.Range("C5").Formula = "=vfbjak"

# Update selection to D5 (this also updates D7 value since C5 is no longer blank)
$ws.Range("D5").Select()
